$d = $word.ActiveDocument

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) Task 1 heading paragraph: split the "List batch mode processing
#    requirements for AdventureWorks" run into a plain-text run + a
#    proofErr-wrapped "AdventureWorks" run.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$xml5 = "<w:p $wns>" + `
    "<w:pPr>" + `
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" + `
      "<w:spacing w:line='285' w:lineRule='atLeast'/>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/><w:b/></w:rPr>" + `
    "</w:pPr>" + `
    "<w:r w:rsidRPr='00DB40B9'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/><w:b/></w:rPr>" + `
      "<w:t xml:space='preserve'>Task 1: </w:t>" + `
    "</w:r>" + `
    "<w:r w:rsidR='00F23397' w:rsidRPr='00F23397'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/><w:b/></w:rPr>" + `
      "<w:t xml:space='preserve'>List batch mode processing requirements for </w:t>" + `
    "</w:r>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/><w:b/></w:rPr>" + `
      "<w:t>AdventureWorks</w:t>" + `
    "</w:r>" + `
    "<w:proofErr w:type='spellEnd'/>" + `
  "</w:p>"
$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------------
# 2) "Use the table below ..." paragraph: drop the _GoBack bookmark and
#    split the trailing " AdventureWorks." run into " " + proofErr-wrapped
#    "AdventureWorks" + ".".
# ---------------------------------------------------------------------------
$p7 = $d.ActiveDocument.Paragraphs(7)
$xml7 = "<w:p $wns>" + `
    "<w:pPr>" + `
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" + `
      "<w:spacing w:line='285' w:lineRule='atLeast'/>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
    "</w:pPr>" + `
    "<w:r w:rsidRPr='009950E6'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t xml:space='preserve'>Use the table below to document </w:t>" + `
    "</w:r>" + `
    "<w:r w:rsidR='00DB40B9' w:rsidRPr='00DB40B9'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t xml:space='preserve'>the requirements that would form part of the </w:t>" + `
    "</w:r>" + `
    "<w:r w:rsidR='00B413BD'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t>b</w:t>" + `
    "</w:r>" + `
    "<w:r w:rsidR='00DB40B9' w:rsidRPr='00DB40B9'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t xml:space='preserve'>atch mode processing of data in an Enterprise BI solution </w:t>" + `
    "</w:r>" + `
    "<w:r w:rsidR='00A4623E'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t>for</w:t>" + `
    "</w:r>" + `
    "<w:r w:rsidR='00DB40B9' w:rsidRPr='00DB40B9'>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t xml:space='preserve'> </w:t>" + `
    "</w:r>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t>AdventureWorks</w:t>" + `
    "</w:r>" + `
    "<w:proofErr w:type='spellEnd'/>" + `
    "<w:r>" + `
      "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
      "<w:t>.</w:t>" + `
    "</w:r>" + `
  "</w:p>"
$p7.Range.InsertXML($xml7)

# ---------------------------------------------------------------------------
# 3) Fill in the table: four requirement rows get descriptive text in the
#    "Requirement" column and a technology suggestion in the third column.
# ---------------------------------------------------------------------------
$t = $d.ActiveDocument.Tables(1)

function Set-CellXml($row, $col, $innerXml) {
    $doc = $word.ActiveDocument
    $tbl = $doc.Tables(1)
    $cell = $tbl.Cell($row, $col)
    $xml = "<w:p $using:wns>" + $innerXml + "</w:p>"
    $cell.Range.InsertXML($xml)
}

# Row 2 (requirement #1): keeps its shd/spacing paragraph formatting.
$inner = "<w:pPr>" + `
    "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" + `
    "<w:spacing w:line='285' w:lineRule='atLeast'/>" + `
    "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
  "</w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t xml:space='preserve'>Ingestion of data from the operational database into </w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>t</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t xml:space='preserve'>he Data Warehouse such as financial reporting. </w:t></w:r>"
Set-CellXml 2 2 $inner

$inner = "<w:pPr><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Azure Data Factory</w:t></w:r>"
Set-CellXml 2 3 $inner

# Row 3 (requirement #2): paragraph formatting gains shd/spacing.
$inner = "<w:pPr>" + `
    "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" + `
    "<w:spacing w:line='285' w:lineRule='atLeast'/>" + `
    "<w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr>" + `
  "</w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Modelling and production of reports for Financial reporting using data in the lake</w:t></w:r>"
Set-CellXml 3 2 $inner

$inner = "<w:pPr><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Azure Synapse</w:t></w:r>"
Set-CellXml 3 3 $inner

# Row 4 (requirement #3): paragraph formatting unchanged (plain).
$inner = "<w:pPr><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Daily Flat Files produced from the Electric Bike Sensors</w:t></w:r>"
Set-CellXml 4 2 $inner

$inner = "<w:pPr><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Databricks</w:t></w:r>"
Set-CellXml 4 3 $inner

# Row 5 (requirement #4): paragraph formatting loses shd/spacing; contains a
# proofErr-wrapped misspelling ("Electrib").
$inner = "<w:pPr><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t xml:space='preserve'>Predictive Maintenance Suggestions from </w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Electrib</w:t></w:r>" + `
  "<w:proofErr w:type='spellEnd'/>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t xml:space='preserve'> Bike sensor data</w:t></w:r>"
Set-CellXml 5 2 $inner

$inner = "<w:pPr><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Segoe UI' w:hAnsi='Segoe UI' w:cs='Segoe UI'/></w:rPr><w:t>Databricks</w:t></w:r>"
Set-CellXml 5 3 $inner
